$d = $word.ActiveDocument

# 1. Remove the first two paragraphs ("Negociações e Vendas:" and its
#    YouTube hyperlink line) - select from the start of paragraph 1
#    through the end of paragraph 2 (including its paragraph mark) and
#    delete the whole range.
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)
$introRange = $d.Range($p1.Range.Start, $p2.Range.End)
$introRange.Delete()

# 2. Split the big run inside the "Aprendizagem Organizacional" hyperlink
#    ("Administração II - Aula 11 - Aprendizagem Organizacional - YouTube")
#    into three runs: "Administ" | "ração II - ... - " | "YouTube".
#    Re-applying (and then clearing) character formatting on a Find match
#    forces Word to materialize a separate run for that match without
#    altering the visible formatting.
$rFirst = $d.Content
$rFirst.Find.Execute("Administ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rFirst.Bold = 1
$rFirst.Bold = 0

$rLast = $d.Content
$rLast.Find.Execute("YouTube", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rLast.Bold = 1
$rLast.Bold = 0

# 3. Turn the trailing empty paragraph into a bookmarked (_GoBack) empty
#    paragraph instead of one holding a stray empty run. Adding a
#    bookmark collapsed at the very end of the document is mishandled by
#    this host, so temporarily extend the document by one character,
#    anchor the bookmark just before it, then remove the scratch
#    character again.
$d.Content.InsertAfter("X")
$lastPara = $d.Paragraphs.Last
$bmRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
$scratch = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$scratch.Delete()
